# Read File -> Get Word Count -> Send to Excel (Date & Word Count)
# Appends the repeated name3/name4 contact block (rows 7-32), a
# date/word-count summary row (33), and a trailing name4 block row (34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-32: repeat the existing 2-row "name3/name4" pattern (rows 5-6) 13 times.
for ($r = 7; $r -le 32; $r += 2) {
    $ws.Cells.Item($r, 1).Value = "name3"
    $ws.Cells.Item($r, 2).Value = "address3"

    $r2 = $r + 1
    $ws.Cells.Item($r2, 1).Value = "name4"
    $ws.Cells.Item($r2, 2).Value = "address4"
    $ws.Cells.Item($r2, 3).Value = "tel4"
    $ws.Cells.Item($r2, 4).Value = "web4"
}

# Row 33: date the file was read + the word count found in it.
# Force text formatting first so the date-like string isn't coerced to a date serial.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "05/02/2022"
$ws.Range("B33").Value = 22

# Row 34: one more trailing "name4" contact block.
$ws.Cells.Item(34, 1).Value = "name4"
$ws.Cells.Item(34, 2).Value = "address4"
$ws.Cells.Item(34, 3).Value = "tel4"
$ws.Cells.Item(34, 4).Value = "web4"
